# Updates the cryptos list (Price / Volume(1h) text columns, plus the
# B/C name+link rotation for rows 44-46) to match the latest scrape,
# per the "Updated cryptos list ... with GitHub Actions" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, $CellRef, $Text)
    $range = $Worksheet.Range($CellRef)
    # Preserve the existing style; force Text number format so values
    # such as "1.001" or "0.3640" are stored as strings (matching the
    # workbook's inline-string cells) instead of being coerced to
    # numbers, then restore the original style/format afterwards.
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $Text
    $range.Style = $origStyle
}

Set-TextValue $ws 'D2' '19.759.99'
Set-TextValue $ws 'E2' '  -8.72%  '
Set-TextValue $ws 'D3' '1.384.46'
Set-TextValue $ws 'E3' '  -9.76%  '
Set-TextValue $ws 'D4' '1.001'
Set-TextValue $ws 'E4' '  +0.10%  '
Set-TextValue $ws 'D5' '1.002'
Set-TextValue $ws 'E5' '  +0.21%  '
Set-TextValue $ws 'D6' '267.07'
Set-TextValue $ws 'E6' '  -7.56%  '
Set-TextValue $ws 'D7' '0.3640'
Set-TextValue $ws 'E7' '  -7.45%  '
Set-TextValue $ws 'D8' '0.3031'
Set-TextValue $ws 'E8' '  -4.46%  '
Set-TextValue $ws 'D9' '38.60'
Set-TextValue $ws 'E9' '  -8.89%  '
Set-TextValue $ws 'D10' '0.9746'
Set-TextValue $ws 'E10' '  -7.58%  '
Set-TextValue $ws 'D11' '0.06386'
Set-TextValue $ws 'E11' '  -10.99%  '
Set-TextValue $ws 'D12' '1.004'
Set-TextValue $ws 'E12' '  +0.33%  '
Set-TextValue $ws 'D13' '5.288'
Set-TextValue $ws 'E13' '  -6.93%  '
Set-TextValue $ws 'D14' '6.048'
Set-TextValue $ws 'E14' '  -8.20%  '
Set-TextValue $ws 'D15' '16.54'
Set-TextValue $ws 'E15' '  -10.87%  '
Set-TextValue $ws 'D16' '1.385.95'
Set-TextValue $ws 'E16' '  -10.09%  '
Set-TextValue $ws 'D17' '0.000009894'
Set-TextValue $ws 'E17' '  -9.50%  '
Set-TextValue $ws 'D18' '0.05603'
Set-TextValue $ws 'E18' '  -14.93%  '
Set-TextValue $ws 'E19' '  +0.27%  '
Set-TextValue $ws 'D20' '69.45'
Set-TextValue $ws 'E20' '  -16.95%  '
Set-TextValue $ws 'D21' '5.488'
Set-TextValue $ws 'E21' '  -10.29%  '
Set-TextValue $ws 'D22' '14.34'
Set-TextValue $ws 'E22' '  -7.18%  '
Set-TextValue $ws 'D23' '10.47'
Set-TextValue $ws 'E23' '  -1.99%  '
Set-TextValue $ws 'D24' '2.251'
Set-TextValue $ws 'E24' '  -4.17%  '
Set-TextValue $ws 'D25' '19.724.49'
Set-TextValue $ws 'E25' '  -8.91%  '
Set-TextValue $ws 'D26' '2.152'
Set-TextValue $ws 'E26' '  -8.41%  '
Set-TextValue $ws 'D27' '135.85'
Set-TextValue $ws 'E27' '  -9.14%  '
Set-TextValue $ws 'E28' '  -10.07%  '
Set-TextValue $ws 'D29' '1.543.81'
Set-TextValue $ws 'E29' '  -9.99%  '
Set-TextValue $ws 'D30' '107.43'
Set-TextValue $ws 'E30' '  -8.20%  '
Set-TextValue $ws 'D31' '3.822'
Set-TextValue $ws 'E31' '  -21.10%  '
Set-TextValue $ws 'D32' '5.199'
Set-TextValue $ws 'E32' '  -14.25%  '
Set-TextValue $ws 'E33' '  -16.00%  '
Set-TextValue $ws 'D34' '0.07573'
Set-TextValue $ws 'E34' '  -6.81%  '
Set-TextValue $ws 'D35' '8.156'
Set-TextValue $ws 'E35' '  -4.21%  '
Set-TextValue $ws 'E36' '  +0.27%  '
Set-TextValue $ws 'D37' '0.05611'
Set-TextValue $ws 'E37' '  -6.74%  '
Set-TextValue $ws 'D38' '4.647'
Set-TextValue $ws 'E38' '  -9.89%  '
Set-TextValue $ws 'D39' '0.02018'
Set-TextValue $ws 'E39' '  -9.31%  '
Set-TextValue $ws 'D40' '0.1864'
Set-TextValue $ws 'E40' '  -7.86%  '
Set-TextValue $ws 'D41' '9.910'
Set-TextValue $ws 'E41' '  -9.13%  '
Set-TextValue $ws 'D42' '1.285'
Set-TextValue $ws 'E42' '  -11.83%  '
Set-TextValue $ws 'D43' '1.042'
Set-TextValue $ws 'E43' '  -11.50%  '
Set-TextValue $ws 'B44' 'EnergySwap'
Set-TextValue $ws 'C44' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws 'D44' '12.00'
Set-TextValue $ws 'E44' '  -7.75%  '
Set-TextValue $ws 'B45' 'TheSandbox'
Set-TextValue $ws 'C45' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws 'D45' '0.5146'
Set-TextValue $ws 'E45' '  -10.72%  '
Set-TextValue $ws 'B46' 'PancakeSwap'
Set-TextValue $ws 'C46' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws 'D46' '3.442'
Set-TextValue $ws 'E46' '  -7.28%  '
Set-TextValue $ws 'D47' '0.4960'
Set-TextValue $ws 'E47' '  -9.88%  '
Set-TextValue $ws 'D48' '108.47'
Set-TextValue $ws 'E48' '  -6.35%  '
Set-TextValue $ws 'D49' '1.710'
Set-TextValue $ws 'E49' '  -8.94%  '
Set-TextValue $ws 'D50' '1.003'
Set-TextValue $ws 'E50' '  +0.35%  '
Set-TextValue $ws 'D51' '1.030'
Set-TextValue $ws 'E51' '  -11.50%  '
